$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in F1, matching the style of the existing header cells (E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Add time_taken values for each data row
$ws.Range("F2").Value = "2021-10-05 13:38:32.193811"
$ws.Range("F3").Value = "2021-10-05 13:38:32.193821"
$ws.Range("F4").Value = "2021-10-05 13:38:32.193824"
$ws.Range("F5").Value = "2021-10-05 13:38:32.193827"
$ws.Range("F6").Value = "2021-10-05 13:38:32.193830"
$ws.Range("F7").Value = "2021-10-05 13:38:32.193832"
